# update scripts wuth new tpm
# Refresh the NATMI LR-pairs output (Mdk-Ptprz1) with recomputed TPM-based
# expression/specificity values for columns G:T across data rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01848533333333334
$ws.Range("N2").Value = 0.05545600000000001
$ws.Range("O2").Value = 0.001625201930372746
$ws.Range("P2").Value = 0.001625201930372746
$ws.Range("Q2").Value = 0.009953335306666667
$ws.Range("R2").Value = 0.08958001776000001
$ws.Range("S2").Value = 0.00005479543833593783
$ws.Range("T2").Value = 0.00005479543833593782
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.002698334581238102
$ws.Range("P3").Value = 0.002698334581238102
$ws.Range("Q3").Value = 0.01652559497666667
$ws.Range("R3").Value = 0.14873035479
$ws.Range("S3").Value = 0.00009097726466645879
$ws.Range("T3").Value = 0.00009097726466645877
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 11.32499966666667
$ws.Range("N4").Value = 33.974999
$ws.Range("O4").Value = 0.9956764634883892
$ws.Range("P4").Value = 0.995676463488389
$ws.Range("Q4").Value = 6.097889445518334
$ws.Range("R4").Value = 54.881005009665
$ws.Range("S4").Value = 0.03357030731874007
$ws.Range("T4").Value = 0.03357030731874006
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.01848533333333334
$ws.Range("N5").Value = 0.05545600000000001
$ws.Range("O5").Value = 0.001625201930372746
$ws.Range("P5").Value = 0.001625201930372746
$ws.Range("Q5").Value = 0.222571355864889
$ws.Range("R5").Value = 2.003142202784
$ws.Range("S5").Value = 0.001225307359782392
$ws.Range("T5").Value = 0.001225307359782391
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.002698334581238102
$ws.Range("P6").Value = 0.002698334581238102
$ws.Range("Q6").Value = 0.3695368403762223
$ws.Range("R6").Value = 3.325831563386001
$ws.Range("S6").Value = 0.002034386718201888
$ws.Range("T6").Value = 0.002034386718201888
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 11.32499966666667
$ws.Range("N7").Value = 33.974999
$ws.Range("O7").Value = 0.9956764634883892
$ws.Range("P7").Value = 0.995676463488389
$ws.Range("Q7").Value = 136.3578619615235
$ws.Range("R7").Value = 1227.220757653711
$ws.Range("S7").Value = 0.7506819158125252
$ws.Range("T7").Value = 0.7506819158125249
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01848533333333334
$ws.Range("N8").Value = 0.05545600000000001
$ws.Range("O8").Value = 0.001625201930372746
$ws.Range("P8").Value = 0.001625201930372746
$ws.Range("Q8").Value = 0.06268564467555555
$ws.Range("R8").Value = 0.56417080208
$ws.Range("S8").Value = 0.0003450991322544165
$ws.Range("T8").Value = 0.0003450991322544164
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.002698334581238102
$ws.Range("P9").Value = 0.002698334581238102
$ws.Range("Q9").Value = 0.1040774316188889
$ws.Range("R9").Value = 0.93669688457
$ws.Range("S9").Value = 0.0005729705983697552
$ws.Range("T9").Value = 0.000572970598369755
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 11.32499966666667
$ws.Range("N10").Value = 33.974999
$ws.Range("O10").Value = 0.9956764634883892
$ws.Range("P10").Value = 0.995676463488389
$ws.Range("Q10").Value = 38.40422524463278
$ws.Range("R10").Value = 345.638027201695
$ws.Range("S10").Value = 0.211424240357124
$ws.Range("T10").Value = 0.211424240357124
